$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rec")

# ---------------------------------------------------------------------------
# 1) Fill in the "checkflag" (M) column for the existing rows 34-37.
#    Column O (realnet = M-D) recalculates automatically.
# ---------------------------------------------------------------------------
$ws.Range("M34").Value = 469
$ws.Range("M35").Value = 166
$ws.Range("M36").Value = 324
$ws.Range("M37").Value = 196

# ---------------------------------------------------------------------------
# 2) Add four new session rows (38-41), copying the "settlement" row format
#    used by row 33 (highlighted fill) and filling in the recorded values.
# ---------------------------------------------------------------------------
$ws.Range("A33:Z33").Copy()
$ws.Range("A38:Z38").PasteSpecial(-4122)
$ws.Range("A39:Z39").PasteSpecial(-4122)
$ws.Range("A40:Z40").PasteSpecial(-4122)
$ws.Range("A41:Z41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# -- row 38 : monoray --------------------------------------------------------
$ws.Range("A38").Value = 42962
$ws.Range("B38").Value = "monoray"
$ws.Range("E38").Value = 3521
$ws.Range("F38").Value = 0.6
$ws.Range("L38").Value = 42964
$ws.Range("Z38").Value = "13486124717"

# -- row 39 : m4n -------------------------------------------------------------
$ws.Range("A39").Value = 42962
$ws.Range("B39").Value = "m4n"
$ws.Range("E39").Value = 5068
$ws.Range("F39").Value = 0.6
$ws.Range("L39").Value = 42964
$ws.Range("Z39").Value = "18657183201"

# -- row 40 : joan --------------------------------------------------------
$ws.Range("A40").Value = 42962
$ws.Range("B40").Value = "joan"
$ws.Range("E40").Value = 3533
$ws.Range("F40").Value = 0.6
$ws.Range("L40").Value = 42964
$ws.Range("Z40").Value = "15824140644"

# -- row 41 : laughing --------------------------------------------------------
$ws.Range("A41").Value = 42962
$ws.Range("B41").Value = "laughing"
$ws.Range("E41").Value = 4177
$ws.Range("F41").Value = 0.61
$ws.Range("L41").Value = 42964
$ws.Range("Z41").Value = "13735467894"

# ---------------------------------------------------------------------------
# 3) Formulas for the new rows, following the same pattern already used by
#    the block above (row 34 = plain formula, rows 35-37 share it).
# ---------------------------------------------------------------------------
$ws.Range("N38").Formula = "=MAX(`$Q`$38:`$Q`$41)/SUBTOTAL(102,`$Q`$38:`$Q`$41)"
$ws.Range("N39:N41").Formula = "=MAX(`$Q`$38:`$Q`$41)/SUBTOTAL(102,`$Q`$38:`$Q`$41)"

$ws.Range("O38:O41").Formula = "=M38-D38"
$ws.Range("P38:P41").Formula = "=IF(K38>D38,INT((K38-D38)*0.95),K38-D38)"
$ws.Range("Q38:Q41").Formula = "=ROUND((L38-A38)*24,2)"
$ws.Range("R38:R41").Formula = "=G38-E38"
$ws.Range("S38:S41").Formula = "=U38-T38"
$ws.Range("T38:T41").Formula = "=INT(E38*F38)"
$ws.Range("U38:U41").Formula = "=INT(G38*H38)"
$ws.Range("V38:V41").Formula = "=S38/R38"
$ws.Range("W38:W41").Formula = "=P38/R38*100/I38"
$ws.Range("X38:X41").Formula = "=P38/Q38/I38"
$ws.Range("Y38:Y41").Formula = "=DATE(YEAR(A38),MONTH(A38),DAY(A38))"

# ---------------------------------------------------------------------------
# 4) Update the view: the frozen pane now keeps the new rows in view and the
#    last-used cell becomes the active selection.
# ---------------------------------------------------------------------------
$ws.Range("G41").Select()
